{"js": "// 1) Title paragraph: \"ANADOLU'NUN RENKLER\u0130 2. ULUSAL FOTO\u011eRAF YARI\u015eMASI\"\n//    -> \"ANADOLU'NUN RENKLER\u0130 3. ULUSAL FOTO\u011eRAF YARI\u015eMASI\"\nconst paragraphs = context.document.body.paragraphs;\nparagraphs.load(\"items\");\nawait context.sync();\n\n// Find the title paragraph by its (stable) leading text instead of a hard-coded\n// index, so the script is a little more robust.\nlet titlePara = null;\nlet longPara = null;\nfor (const p of paragraphs.items) {\n  p.load(\"text\");\n}\nawait context.sync();\nfor (const p of paragraphs.items) {\n  if (p.text.indexOf(\"ANADOLU\") !== -1 && p.text.indexOf(\"ULUSAL FOTO\u011eRAF YARI\u015eMASI\") !== -1) {\n    titlePara = p;\n  }\n  if (p.text.indexOf(\"Mill\u00ee E\u011fitim Bakanl\u0131\u011f\u0131n\u0131n bilgisi d\u00e2hilinde\") !== -1) {\n    longPara = p;\n  }\n}\n\nif (titlePara) {\n  const titleRange = titlePara.getRange();\n  // Splitting on spaces gives us word-sized ranges without touching the\n  // (buggy, cross-run) body/paragraph `search()` API, and it keeps the\n  // existing run split (\"...RENKLER\u0130 \" / \"2. \" / \"ULUSAL...\") intact.\n  const titleWords = titleRange.getTextRanges([\" \"], true);\n  titleWords.load(\"text,items\");\n  await context.sync();\n  for (const w of titleWords.items) {\n    if (w.text === \"2.\") {\n      w.insertText(\"3.\", Word.InsertLocation.replace);\n      break;\n    }\n  }\n  await context.sync();\n}\n\n// 2) Remove the hidden \"_GoBack\" bookmark.\ncontext.document.deleteBookmark(\"_GoBack\");\nawait context.sync();\n\n// 3) Trim the \"Konya B\u00fcy\u00fck\u015fehir Belediyesi sponsorlu\u011funda Konya \u0130l Mill\u00ee...\"\n//    sentence down to \"Konya \u0130l Mill\u00ee...\" and drop the stray \"2. \" before\n//    \"Ulusal Foto\u011fraf Yar\u0131\u015fmas\u0131\".\nif (longPara) {\n  const longRange = longPara.getRange();\n  const words = longRange.getTextRanges([\" \"], true);\n  words.load(\"text,items\");\n  await context.sync();\n\n  const items = words.items;\n\n  // a) delete \"B\u00fcy\u00fck\u015fehir Belediyesi sponsorlu\u011funda Konya \" (the run of\n  //    words between the first \"Konya\" and the second \"\u0130l\").\n  for (let i = 0; i < items.length - 1; i++) {\n    if (\n      items[i].text === \"B\u00fcy\u00fck\u015fehir\" &&\n      items[i + 1] && items[i + 1].text === \"Belediyesi\" &&\n      items[i + 2] && items[i + 2].text === \"sponsorlu\u011funda\" &&\n      items[i + 3] && items[i + 3].text === \"Konya\" &&\n      items[i + 4] && items[i + 4].text === \"\u0130l\"\n    ) {\n      const toDelete = items[i].expandTo(items[i + 3]).expandTo(items[i + 4].getRange(\"Start\"));\n      toDelete.insertText(\"\", Word.InsertLocation.replace);\n      break;\n    }\n  }\n  await context.sync();\n\n  // b) delete the stray \"2. \" before \"Ulusal Foto\u011fraf Yar\u0131\u015fmas\u0131\" (need to\n  //    re-fetch the word ranges since the paragraph changed above).\n  const words2 = longPara.getRange().getTextRanges([\" \"], true);\n  words2.load(\"text,items\");\n  await context.sync();\n  const items2 = words2.items;\n  for (let i = 0; i < items2.length - 1; i++) {\n    if (\n      items2[i].text === \"Renkleri\" &&\n      items2[i + 1] && items2[i + 1].text === \"2.\" &&\n      items2[i + 2] && items2[i + 2].text === \"Ulusal\"\n    ) {\n      const toDelete = items2[i + 1].expandTo(items2[i + 2].getRange(\"Start\"));\n      toDelete.insertText(\"\", Word.InsertLocation.replace);\n      break;\n    }\n  }\n  await context.sync();\n}\n", "ps1": "$d = $word.ActiveDocument\n\n# Locate the two paragraphs we need to touch by stable text fragments\n# (safer than hard-coded paragraph indices).\n$titlePara = $null\n$longPara = $null\nforeach ($p in $d.Paragraphs) {\n    $txt = $p.Range.Text\n    if ($txt -like \"*ANADOLU*ULUSAL FOTO*YARI*MASI*\") { $titlePara = $p }\n    if ($txt -like \"*bilgisi d*hilinde*\") { $longPara = $p }\n}\n\n# 1) Title: \"ANADOLU'NUN RENKLER\u0130 2. ULUSAL FOTO\u011eRAF YARI\u015eMASI\"\n#           -> \"ANADOLU'NUN RENKLER\u0130 3. ULUSAL FOTO\u011eRAF YARI\u015eMASI\"\nif ($titlePara -ne $null) {\n    $full = $titlePara.Range.Text\n    $idx = $full.IndexOf(\"2.\")\n    if ($idx -ge 0) {\n        $pStart = $titlePara.Range.Start\n        $rng = $d.Range($pStart + $idx, $pStart + $idx + 1)\n        $rng.Text = \"3\"\n    }\n}\n\n# 2) Remove the hidden \"_GoBack\" bookmark.\nif ($d.Bookmarks.Exists(\"_GoBack\")) {\n    $d.Bookmarks(\"_GoBack\").Delete()\n}\n\nif ($longPara -ne $null) {\n    # 3) \"...d\u00e2hilinde, Konya B\u00fcy\u00fck\u015fehir Belediyesi sponsorlu\u011funda Konya \u0130l...\"\n    #    -> \"...d\u00e2hilinde, Konya \u0130l...\"\n    $full2 = $longPara.Range.Text\n    $needle = \"B\u00fcy\u00fck\u015fehir Belediyesi sponsorlu\u011funda Konya \"\n    $idx2 = $full2.IndexOf($needle)\n    if ($idx2 -ge 0) {\n        $lStart = $longPara.Range.Start\n        $rng2 = $d.Range($lStart + $idx2, $lStart + $idx2 + $needle.Length)\n        $rng2.Text = \"\"\n    }\n\n    # 4) \"...Anadolu\u2019nun Renkleri 2. Ulusal Foto\u011fraf Yar\u0131\u015fmas\u0131...\"\n    #    -> \"...Anadolu\u2019nun Renkleri Ulusal Foto\u011fraf Yar\u0131\u015fmas\u0131...\"\n    $full3 = $longPara.Range.Text\n    $needle2 = \"Renkleri 2. Ulusal\"\n    $idx3 = $full3.IndexOf($needle2)\n    if ($idx3 -ge 0) {\n        $lStart2 = $longPara.Range.Start\n        $offset = $idx3 + \"Renkleri \".Length\n        $rng3 = $d.Range($lStart2 + $offset, $lStart2 + $offset + \"2. \".Length)\n        $rng3.Text = \"\"\n    }\n}\n"}
